$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.719.93"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "2.524.94"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("D4").Formula = "=""0.999"""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Formula = "=""309.72"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").Formula = "=""101.59"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +3.87%  "
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Formula = "=""0.526"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -2.50%  "
$ws.Range("D10").Formula = "=""35.95"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").Formula = "=""0.0803"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Formula = "=""7.30"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -3.16%  "
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "2.912.44"
$ws.Range("E14").Value = "  +3.71%  "
$ws.Range("D15").Formula = "=""15.63"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").Value = "2.529.03"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Formula = "=""0.806"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -5.14%  "
$ws.Range("D18").Value = "42.688.34"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Formula = "=""6.73"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").Value = "0.0₃0949"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Formula = "=""12.15"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -5.15%  "
$ws.Range("D22").Formula = "=""69.24"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Formula = "=""244.07"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -3.94%  "
$ws.Range("E24").Value = "  -3.10%  "
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Formula = "=""26.15"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -4.44%  "
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("D29").Formula = "=""39.14"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -5.06%  "
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("D31").Formula = "=""156.46"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Formula = "=""5.77"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").Formula = "=""2.78"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +10.44%  "
$ws.Range("D34").Formula = "=""0.0786"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -3.10%  "
$ws.Range("D35").Formula = "=""2.62"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("E36").Value = "  -6.54%  "
$ws.Range("D37").Formula = "=""3.19"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -7.85%  "
$ws.Range("D38").Formula = "=""18.22"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("E41").Value = "  +5.83%  "
$ws.Range("D42").Formula = "=""21.94"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("D45").Formula = "=""3.28"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("D46").Value = "1.990.61"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Formula = "=""8.83"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").Value = "2.766.29"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Formula = "=""80.20"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -4.01%  "
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("D51").Formula = "=""72.19"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -3.42%  "
$excel.CutCopyMode = 0
